# Update the "check version" timestamp stored in column J.
# Every cell in column J that still carries the old capture timestamp
# (1587083823.630226) is bumped to the new one (1587089010.639998).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = 1587083823.630226
$newValue = 1587089010.639998

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # column J = 10
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
